$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Step 1: convert A135:A143 from text "NN" to numeric NN (unchanged value) ---
foreach ($r in 135..143) {
    $n = [double]$ws.Cells.Item($r, 1).Text
    $ws.Cells.Item($r, 1).Value = $n
}

# --- Step 2: append new rows 144-148 ---
# Row 144
$ws.Cells.Item(144, 1).NumberFormat = "@"
$ws.Cells.Item(144, 1).Value = '27'
$ws.Cells.Item(144, 2).Value = ' saas demo video - Upwork'
$ws.Cells.Item(144, 3).Value = 'https://www.upwork.com/jobs/saas-demo-video_%7E01f22de24fc078d2dc?source=rss'
$ws.Cells.Item(144, 4).Value = 'I need a saas demo video for my software and I want it about 30-40 seconds long I''ll provide the link of figma from there you can pick up screens. No voice over is needed. Just simple animations.
Budget
: $45
Posted On
: June 15, 2024 11:08 UTC
Category
: Video Production
Skills
:Video Editing,     Explainer Video,     Video Post-Editing,     Video Production,     Motion Graphics    
Skills
:        Video Editing,                     Explainer Video,                     Video Post-Editing,                     Video Production,                     Motion Graphics            
Country
: United States
click to apply
'
$ws.Cells.Item(144, 5).Value = 'I need a saas demo video for my software and I want it about 30-40 seconds long I&#039;ll provide the link of figma from there you can pick up screens. No voice over is needed. Just simple&nbsp;animations.<br /><br /><b>Budget</b>: $45
<br /><b>Posted On</b>: June 15, 2024 11:08 UTC<br /><b>Category</b>: Video Production<br /><b>Skills</b>:Video Editing,     Explainer Video,     Video Post-Editing,     Video Production,     Motion Graphics    
<br /><b>Skills</b>:        Video Editing,                     Explainer Video,                     Video Post-Editing,                     Video Production,                     Motion Graphics            <br /><b>Country</b>: United States
<br /><a href="https://www.upwork.com/jobs/saas-demo-video_%7E01f22de24fc078d2dc?source=rss">click to apply</a>
'
$ws.Cells.Item(144, 6).Value = 'Sat, 15 Jun 2024 11:08:06 +0000'
$ws.Cells.Item(144, 7).Value = 'https://www.upwork.com/jobs/saas-demo-video_%7E01f22de24fc078d2dc?source=rss'
$ws.Cells.Item(144, 9).NumberFormat = "@"
$ws.Cells.Item(144, 9).Value = '$45'
$ws.Cells.Item(144, 10).Value = 'June 15, 2024 11:08 UTC'
$ws.Cells.Item(144, 11).Value = 'Video Production'
$ws.Cells.Item(144, 12).Value = 'Video Editing,     Explainer Video,     Video Post-Editing,     Video Production,     Motion Graphics'
$ws.Cells.Item(144, 13).Value = 'United States'

# Row 145
$ws.Cells.Item(145, 1).NumberFormat = "@"
$ws.Cells.Item(145, 1).Value = '27'
$ws.Cells.Item(145, 2).Value = 'TikTok Video Editor (Monthly)  - Upwork'
$ws.Cells.Item(145, 3).Value = 'https://www.upwork.com/jobs/TikTok-Video-Editor-Monthly_%7E018319ac8ae6b5cd0c?source=rss'
$ws.Cells.Item(145, 4).Value = 'Hello!
I am looking for a freelancer that will help me edit simple videos for TikTok. 
The videos are as said, simple. Most of the only require cutting and captions. I will share some more examples after the NDA has been signed. (Do not apply if you''re not ready to sign the non-disclosure agreement)
But here are some close references of what kind of content we are looking for - (https://www.tiktok.com/@musicmediaco/video/7201943783473155370?q=music%20media&amp;t=1717716126194) 
This project is a long-term project and we will award the right freelancer with a long-term collaboration.
The MOST IMPORTANT qualities that we are looking for are:
1. Will to learn and the ability to adapt
2. Prompt communication
3. Adhere to set deadlines 
Please send a link to your portfolio, attach your CV and also tell me a fun fact about your favorite animal, so I know you''ve read this. 
Thank you very much, I''m looking forward to your applications. 
I will check every single one of the applications. 
-D
Budget
: $200
Posted On
: June 15, 2024 11:05 UTC
Category
: Video Editing
Skills
:Audio Editing,     Compositing,     DaVinci Resolve,     Video Editing,     Adobe After Effects,     Adobe Premiere Pro,     Video Post-Editing,     Video Production,     Music Video    
Skills
:        Audio Editing,                     Compositing,                     DaVinci Resolve,                     Video Editing,                     Adobe After Effects,                     Adobe Premiere Pro,                     Video Post-Editing,                     Video Production,                     Music Video            
Country
: Croatia
click to apply
'
$ws.Cells.Item(145, 5).Value = 'Hello!<br /><br />
I am looking for a freelancer that will help me edit simple videos for TikTok. <br /><br />
The videos are as said, simple. Most of the only require cutting and captions. I will share some more examples after the NDA has been signed. (Do not apply if you&#039;re not ready to sign the non-disclosure agreement)<br /><br />
But here are some close references of what kind of content we are looking for - (https://www.tiktok.com/@musicmediaco/video/7201943783473155370?q=music%20media&amp;amp;t=1717716126194) <br /><br />
This project is a long-term project and we will award the right freelancer with a long-term collaboration.<br /><br />
The MOST IMPORTANT qualities that we are looking for are:<br />
1. Will to learn and the ability to adapt<br />
2. Prompt communication<br />
3. Adhere to set deadlines <br /><br />
Please send a link to your portfolio, attach your CV and also tell me a fun fact about your favorite animal, so I know you&#039;ve read this. <br /><br />
Thank you very much, I&#039;m looking forward to your applications. <br />
I will check every single one of the applications. <br /><br />
-D<br /><br /><b>Budget</b>: $200
<br /><b>Posted On</b>: June 15, 2024 11:05 UTC<br /><b>Category</b>: Video Editing<br /><b>Skills</b>:Audio Editing,     Compositing,     DaVinci Resolve,     Video Editing,     Adobe After Effects,     Adobe Premiere Pro,     Video Post-Editing,     Video Production,     Music Video    
<br /><b>Skills</b>:        Audio Editing,                     Compositing,                     DaVinci Resolve,                     Video Editing,                     Adobe After Effects,                     Adobe Premiere Pro,                     Video Post-Editing,                     Video Production,                     Music Video            <br /><b>Country</b>: Croatia
<br /><a href="https://www.upwork.com/jobs/TikTok-Video-Editor-Monthly_%7E018319ac8ae6b5cd0c?source=rss">click to apply</a>
'
$ws.Cells.Item(145, 6).Value = 'Sat, 15 Jun 2024 11:05:56 +0000'
$ws.Cells.Item(145, 7).Value = 'https://www.upwork.com/jobs/TikTok-Video-Editor-Monthly_%7E018319ac8ae6b5cd0c?source=rss'
$ws.Cells.Item(145, 9).NumberFormat = "@"
$ws.Cells.Item(145, 9).Value = '$200'
$ws.Cells.Item(145, 10).Value = 'June 15, 2024 11:05 UTC'
$ws.Cells.Item(145, 11).Value = 'Video Editing'
$ws.Cells.Item(145, 12).Value = 'Audio Editing,     Compositing,     DaVinci Resolve,     Video Editing,     Adobe After Effects,     Adobe Premiere Pro,     Video Post-Editing,     Video Production,     Music Video'
$ws.Cells.Item(145, 13).Value = 'Croatia'

# Row 146
$ws.Cells.Item(146, 1).NumberFormat = "@"
$ws.Cells.Item(146, 1).Value = '27'
$ws.Cells.Item(146, 2).Value = 'Shopify Brand Product Photo/Tile Creation - Upwork'
$ws.Cells.Item(146, 3).Value = 'https://www.upwork.com/jobs/Shopify-Brand-Product-Photo-Tile-Creation_%7E01aea4b8be925905e0?source=rss'
$ws.Cells.Item(146, 4).Value = 'Hello,
I am looking for someone who can help me create the best possible product pictures for my brand. We just did a big photo shoot and had a lot of UGC, 3D models of our products, and videos to use. Specifically Jewelry.
I am looking for someone who will use all of our existing content, to create the best pictures that will sell the product.
Main Example 1: https://moonmagic.com/collections/all-bestsellers/products/flow-ring-stardust-band?variant=39469970784328
Example 2: https://spacegoods.com/en-eu/products/rainbow-dust-premium-starter-kit-coffee
Example 3: https://podcompany.com/products/the-ice-pod
As you can see, they all use a mix of photos of the actual products, renders with USPs, vibe photos, graphic design on a photo, etc.
We need to do the same for our jewelry brand, with all of our existing content, call out all the USPs in the product photos, select the best photos, use videos, really sell the product, all in the style and vibe of our brand, and new website.
Please only apply if you have previous experience and include previous work so I can see what you have done.
Greetings,
Tristan
Hourly Range
: $6.00-$16.00
Posted On
: June 15, 2024 10:45 UTC
Category
: Image Editing
Skills
:Graphic Design,     Adobe Photoshop,     Photo Editing    
Skills
:        Graphic Design,                     Adobe Photoshop,                     Photo Editing            
Country
: Netherlands
click to apply
'
$ws.Cells.Item(146, 5).Value = 'Hello,<br /><br />
I am looking for someone who can help me create the best possible product pictures for my brand. We just did a big photo shoot and had a lot of UGC, 3D models of our products, and videos to use. Specifically Jewelry.<br /><br />
I am looking for someone who will use all of our existing content, to create the best pictures that will sell the product.<br /><br />
Main Example 1: https://moonmagic.com/collections/all-bestsellers/products/flow-ring-stardust-band?variant=39469970784328<br /><br />
Example 2: https://spacegoods.com/en-eu/products/rainbow-dust-premium-starter-kit-coffee<br />
Example 3: https://podcompany.com/products/the-ice-pod<br /><br />
As you can see, they all use a mix of photos of the actual products, renders with USPs, vibe photos, graphic design on a photo, etc.<br /><br />
We need to do the same for our jewelry brand, with all of our existing content, call out all the USPs in the product photos, select the best photos, use videos, really sell the product, all in the style and vibe of our brand, and new website.<br /><br />
Please only apply if you have previous experience and include previous work so I can see what you have done.<br /><br />
Greetings,<br />
Tristan<br /><br /><br /><b>Hourly Range</b>: $6.00-$16.00
<br /><b>Posted On</b>: June 15, 2024 10:45 UTC<br /><b>Category</b>: Image Editing<br /><b>Skills</b>:Graphic Design,     Adobe Photoshop,     Photo Editing    
<br /><b>Skills</b>:        Graphic Design,                     Adobe Photoshop,                     Photo Editing            <br /><b>Country</b>: Netherlands
<br /><a href="https://www.upwork.com/jobs/Shopify-Brand-Product-Photo-Tile-Creation_%7E01aea4b8be925905e0?source=rss">click to apply</a>
'
$ws.Cells.Item(146, 6).Value = 'Sat, 15 Jun 2024 10:45:56 +0000'
$ws.Cells.Item(146, 7).Value = 'https://www.upwork.com/jobs/Shopify-Brand-Product-Photo-Tile-Creation_%7E01aea4b8be925905e0?source=rss'
$ws.Cells.Item(146, 8).Value = '$6.00-$16.00'
$ws.Cells.Item(146, 10).Value = 'June 15, 2024 10:45 UTC'
$ws.Cells.Item(146, 11).Value = 'Image Editing'
$ws.Cells.Item(146, 12).Value = 'Graphic Design,     Adobe Photoshop,     Photo Editing'
$ws.Cells.Item(146, 13).Value = 'Netherlands'

# Row 147
$ws.Cells.Item(147, 1).NumberFormat = "@"
$ws.Cells.Item(147, 1).Value = '27'
$ws.Cells.Item(147, 2).Value = 'Video Editor for Youtube Channel - Glow Up Niche - Upwork'
$ws.Cells.Item(147, 3).Value = 'https://www.upwork.com/jobs/Video-Editor-for-Youtube-Channel-Glow-Niche_%7E013878b70714d31602?source=rss'
$ws.Cells.Item(147, 4).Value = 'We are looking for a skilled video editor to join our team and work on our YouTube channel in the glow up niche. The ideal candidate will have experience editing videos for YouTube with a strong focus on creating visually appealing and engaging content. The main responsibilities will include editing raw footage, adding transitions and effects, optimizing video and audio quality, and ensuring the final product aligns with our brand''s aesthetic. The video editor should have a creative mindset and be able to follow guidelines to create consistent content. Proficiency with video editing software (e.g., Adobe Premiere Pro, Final Cut Pro) and knowledge of YouTube''s best practices are required. Need to edit exactly like this: (link removed) , Paying 25 dollars per 5 minute video
Budget
: $500
Posted On
: June 15, 2024 10:16 UTC
Category
: Video Editing
Skills
:Video Editing,     Video Post-Editing,     Adobe Premiere Pro,     Video Production    
Skills
:        Video Editing,                     Video Post-Editing,                     Adobe Premiere Pro,                     Video Production            
Country
: Czech Republic
click to apply
'
$ws.Cells.Item(147, 5).Value = 'We are looking for a skilled video editor to join our team and work on our YouTube channel in the glow up niche. The ideal candidate will have experience editing videos for YouTube with a strong focus on creating visually appealing and engaging content. The main responsibilities will include editing raw footage, adding transitions and effects, optimizing video and audio quality, and ensuring the final product aligns with our brand&#039;s aesthetic. The video editor should have a creative mindset and be able to follow guidelines to create consistent content. Proficiency with video editing software (e.g., Adobe Premiere Pro, Final Cut Pro) and knowledge of YouTube&#039;s best practices are required. Need to edit exactly like this: (link removed) , Paying 25 dollars per 5 minute video<br /><br /><b>Budget</b>: $500
<br /><b>Posted On</b>: June 15, 2024 10:16 UTC<br /><b>Category</b>: Video Editing<br /><b>Skills</b>:Video Editing,     Video Post-Editing,     Adobe Premiere Pro,     Video Production    
<br /><b>Skills</b>:        Video Editing,                     Video Post-Editing,                     Adobe Premiere Pro,                     Video Production            <br /><b>Country</b>: Czech Republic
<br /><a href="https://www.upwork.com/jobs/Video-Editor-for-Youtube-Channel-Glow-Niche_%7E013878b70714d31602?source=rss">click to apply</a>
'
$ws.Cells.Item(147, 6).Value = 'Sat, 15 Jun 2024 10:16:15 +0000'
$ws.Cells.Item(147, 7).Value = 'https://www.upwork.com/jobs/Video-Editor-for-Youtube-Channel-Glow-Niche_%7E013878b70714d31602?source=rss'
$ws.Cells.Item(147, 9).NumberFormat = "@"
$ws.Cells.Item(147, 9).Value = '$500'
$ws.Cells.Item(147, 10).Value = 'June 15, 2024 10:16 UTC'
$ws.Cells.Item(147, 11).Value = 'Video Editing'
$ws.Cells.Item(147, 12).Value = 'Video Editing,     Video Post-Editing,     Adobe Premiere Pro,     Video Production'
$ws.Cells.Item(147, 13).Value = 'Czech Republic'

# Row 148
$ws.Cells.Item(148, 1).NumberFormat = "@"
$ws.Cells.Item(148, 1).Value = '27'
$ws.Cells.Item(148, 2).Value = 'Looking for UFC Video CREATOR for LONG TERM WORK - Upwork'
$ws.Cells.Item(148, 3).Value = 'https://www.upwork.com/jobs/Looking-for-UFC-Video-CREATOR-for-LONG-TERM-WORK_%7E01f6b973bfba1ed838?source=rss'
$ws.Cells.Item(148, 4).Value = 'We are seeking a talented and experienced Video Editor to join our YouTube channel''s editing team. As a Video Editor, you will play a crucial role in producing high-quality videos centered around the exciting world of UFC and MMA. If you have a passion for combat sports and possess exceptional video editing skills, this opportunity is perfect for you.
The channel is a MMA/UFC news channel so when there is a topic to post about, we have to work quick. 
Requirments: 
Knowledge about the UFC and MMA world 
Edit 2-3 videos per week 
We have a style of videos we make and want you to make the same. 
Note: 
We only look for people who have knowledge about the UFC world. If you are not, this position is not for you. 
Hourly Range
: $20.00-$35.00
Posted On
: June 15, 2024 10:16 UTC
Category
: Video Editing
Skills
:Video Editing    
Skills
:        Video Editing            
Country
: Georgia
click to apply
'
$ws.Cells.Item(148, 5).Value = 'We are seeking a talented and experienced Video Editor to join our YouTube channel&#039;s editing team. As a Video Editor, you will play a crucial role in producing high-quality videos centered around the exciting world of UFC and MMA. If you have a passion for combat sports and possess exceptional video editing skills, this opportunity is perfect for you.<br /><br />
The channel is a MMA/UFC news channel so when there is a topic to post about, we have to work quick. <br /><br />
Requirments: <br />
Knowledge about the UFC and MMA world <br />
Edit 2-3 videos per week <br /><br /><br />
We have a style of videos we make and want you to make the same. <br /><br />
Note: <br />
We only look for people who have knowledge about the UFC world. If you are not, this position is not for you. <br /><br /><br /><b>Hourly Range</b>: $20.00-$35.00
<br /><b>Posted On</b>: June 15, 2024 10:16 UTC<br /><b>Category</b>: Video Editing<br /><b>Skills</b>:Video Editing    
<br /><b>Skills</b>:        Video Editing            <br /><b>Country</b>: Georgia
<br /><a href="https://www.upwork.com/jobs/Looking-for-UFC-Video-CREATOR-for-LONG-TERM-WORK_%7E01f6b973bfba1ed838?source=rss">click to apply</a>
'
$ws.Cells.Item(148, 6).Value = 'Sat, 15 Jun 2024 10:16:15 +0000'
$ws.Cells.Item(148, 7).Value = 'https://www.upwork.com/jobs/Looking-for-UFC-Video-CREATOR-for-LONG-TERM-WORK_%7E01f6b973bfba1ed838?source=rss'
$ws.Cells.Item(148, 8).Value = '$20.00-$35.00'
$ws.Cells.Item(148, 10).Value = 'June 15, 2024 10:16 UTC'
$ws.Cells.Item(148, 11).Value = 'Video Editing'
$ws.Cells.Item(148, 12).Value = 'Video Editing'
$ws.Cells.Item(148, 13).Value = 'Georgia'

